$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 3
$ws.Range("B5").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 5
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 2
$ws.Range("B18").Value = 4
$ws.Range("B20").Value = 2
$ws.Range("B22").Value = 1
$ws.Range("B26").Value = 1
$ws.Range("B27").Value = 5
$ws.Range("B28").Value = 1
$ws.Range("B29").Value = 1
$ws.Range("B30").Value = 4
$ws.Range("B31").Value = 5
$ws.Range("B32").Value = 4
$ws.Range("B33").Value = 5
$ws.Range("B34").Value = 3
$ws.Range("B38").Value = 4
$ws.Range("B39").Value = 4
$ws.Range("B40").Value = 5
$ws.Range("B41").Value = 1
$ws.Range("B42").Value = 2
$ws.Range("B44").Value = 4
$ws.Range("B45").Value = 3
$ws.Range("B46").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("B48").Value = 3
$ws.Range("B49").Value = 3
$ws.Range("B52").Value = 4
$ws.Range("B53").Value = 2
$ws.Range("B55").Value = 2
$ws.Range("B56").Value = 2
$ws.Range("B58").Value = 4
$ws.Range("B60").Value = 4
$ws.Range("B61").Value = 3
$ws.Range("B63").Value = 3
$ws.Range("B64").Value = 4
$ws.Range("B65").Value = 3
$ws.Range("B69").Value = 1
$ws.Range("B70").Value = 4
$ws.Range("B71").Value = 1
$ws.Range("B72").Value = 3
$ws.Range("B74").Value = 4
$ws.Range("B75").Value = 2
$ws.Range("B77").Value = 2
$ws.Range("B78").Value = 2
$ws.Range("B80").Value = 3
$ws.Range("B81").Value = 2
$ws.Range("B83").Value = 1
$ws.Range("B84").Value = 2
$ws.Range("B87").Value = 5
$ws.Range("B88").Value = 5
$ws.Range("B89").Value = 3
$ws.Range("B95").Value = 4
$ws.Range("B96").Value = 1
$ws.Range("B97").Value = 1
$ws.Range("B98").Value = 3
$ws.Range("B99").Value = 3
$ws.Range("B100").Value = 1
$ws.Range("B101").Value = 1
